$p = $ppt.ActivePresentation
Write-Host "Slide count before: $($p.Slides.Count)"
$s = $p.Slides.Add($p.Slides.Count + 1, 12)
Write-Host "Slide count after: $($p.Slides.Count)"
